$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers.
$ws.Range("Q2").Value = 663735
$ws.Range("R2").Value = 6710544

# Clear the "Starttid" / "Sluttid" (start/end time) cells — times unknown.
$ws.Range("Z2").Value = ""
$ws.Range("AB2").Value = ""
